$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-52.9(se=7.9)"
$ws.Range("B3").Value = "-24.3(se=2.8)"
$ws.Range("C3").Value = "0.28(95% CI, 0.095-0.42)"
$ws.Range("F3").Value = "1.4(95% CI, 0.48-2)"
$ws.Range("B4").Value = "-34.4(se=4.5)"
$ws.Range("B5").Value = "-59(se=18.2)"
$ws.Range("F6").Value = "1.3(95% CI, 0.69-2.1)"
